# Update "想去人数" (want-to-go count) values in F column on the
# "展览" sheet and the mirrored "全部类型" sheet, per gh-pages data
# regeneration at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll = $wb.Worksheets.Item("全部类型")

# row -> new F-column value, for the "展览" sheet
$exhibitUpdates = @{
    4  = 75
    5  = 1665
    6  = 3266
    7  = 829
    8  = 2063
    9  = 1975
    10 = 1020
    11 = 354
    13 = 1614
    18 = 75
    19 = 1447
    20 = 531
    21 = 636
    22 = 323
    23 = 10754
    24 = 11736
    25 = 857
    29 = 448
}

# row -> new F-column value, for the "全部类型" sheet (same events,
# offset rows since this sheet concatenates multiple categories)
$allUpdates = @{
    6  = 75
    7  = 1665
    8  = 3266
    9  = 829
    10 = 2063
    11 = 1975
    12 = 1020
    13 = 354
    15 = 1614
    22 = 75
    23 = 1447
    24 = 531
    25 = 636
    26 = 323
    27 = 10754
    28 = 11736
    29 = 857
    35 = 448
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
